# Edit script: fix UI and translate to vietnamese
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data values (order matches shared-string insertion order in target)
$ws.Range("E2").Value = "Đo thiết bị"
$ws.Range("G2").Value = " mg/l"
$ws.Range("A2").Value = "f7067e30-4bb1-4812-99fb-db4639e6ca04"
$ws.Range("B2").Value = "NO"

# Update the selected cell on the sheet (UI fix)
$ws.Range("D15").Select()
